$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update existing entry's label: "Words from books" -> "Words from books and newspapers"
#    (row 31: A31=54)
$ws.Range("B31").Value = "Words from books and newspapers"

# 2. Renumber "Finding a place to live, apartments, renting" from 70 -> 59 (row 43)
$ws.Range("A43").Value = 59

# 3. Renumber "Past tense" from 82 -> 55 (row 54)
$ws.Range("A54").Value = 55

# 4. Append two brand-new rows at the bottom of the list
$ws.Range("A70").Value = 88
$ws.Range("B70").Value = "Church"

$ws.Range("A71").Value = 80
$ws.Range("B71").Value = "Miscellaneous 6"

# 5. Re-sort the whole list (A1:B71) ascending by column A, mirroring the
#    worksheet's existing sortState/sortCondition behavior
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A1:A71"))
$sortObj.SetRange($ws.Range("A1:B71"))
$sortObj.Header = 2
[void]$sortObj.Apply()

# 6. Update the sheet view to match the final cursor/scroll position
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B26").Select()
